# CKM Variable List.xlsx - add "interview_period" (RIDEXMON) variable
#
# The commit adds a new variable row to the DEMO sheet's variable list:
#   interview_period -> RIDEXMON (repeated across every NHANES cycle column)
# inserted right after the "age"/RIDAGEYR row (i.e. as the new row 4), pushing
# every row below it down by one. It also leaves the workbook positioned with
# the DEMO sheet active/selected (as it was the sheet last edited), while the
# previously-active DIQ sheet and the INQ sheet keep the cursor positions they
# were left at.

$wb = $excel.ActiveWorkbook

# --- DIQ sheet: was the active tab before the edit; record where its cursor
#     ended up once it was no longer the active sheet. ---
$diq = $wb.Worksheets.Item("DIQ")
$diq.Activate()
$diq.Range("A30").Select()

# --- INQ sheet: cursor moved while reviewing the workbook. ---
$inq = $wb.Worksheets.Item("INQ")
$inq.Activate()
$inq.Range("E16").Select()

# --- DEMO sheet: the actual content edit. ---
$demo = $wb.Worksheets.Item("DEMO")
$demo.Activate()

# Insert a new row above the current row 4 ("gender"/RIAGENDR), shifting the
# remainder of the variable list down by one row.
$demo.Rows.Item(4).Insert()

# Column A holds the friendly variable name, columns B:M hold the raw NHANES
# field code for each survey cycle (1999-2000 ... 2021-2023).
$demo.Range("A4").Value = "interview_period"
$demo.Range("B4:M4").Value = "RIDEXMON"

# Leave the cursor/active sheet where the author left it after the edit.
$demo.Range("F27").Select()
